$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, shifting existing rows 46-69 down to 47-70.
$ws.Rows(46).Insert()

# Populate the newly inserted row 46 with the new weekly price record.
$ws.Cells.Item(46,1).Value  = 1
$ws.Cells.Item(46,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(46,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(46,4).Value  = Get-Date -Year 2022 -Month 6 -Day 14 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(46,5).Value  = 15
$ws.Cells.Item(46,6).Value  = 100112012
$ws.Cells.Item(46,7).Value  = "Espinaca"
$ws.Cells.Item(46,8).Value  = "Sin especificar"
$ws.Cells.Item(46,9).Value  = "Primera"
$ws.Cells.Item(46,10).Value = 300
$ws.Cells.Item(46,11).Value = 2800
$ws.Cells.Item(46,12).Value = 3000
$ws.Cells.Item(46,13).Value = 2900
$ws.Cells.Item(46,14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(46,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(46,16).Value = 967
$ws.Cells.Item(46,17).Value = 3
$ws.Cells.Item(46,18).Value = "Hortaliza"
